# Applies the "Added one hot encoding" change to the Codebook sheet:
#   1. Fixes the I14 sample-values list (drops 'OTHER:GERMAN_STUDIES', adds 'GRS').
#   2. Widens column A from 23 to 25 characters.
#   3. Appends 13 new one-hot-encoded variable rows (70-82) describing the
#      Media_* social-media CE items and the News_* local-news-source items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codebook")

# --- 1. Fix existing sample values in I14 -------------------------------
$ws.Range("I14").Value = "['THE', 'CMS', 'GRS', 'EDU', 'ANT', 'PSY', 'ENV', 'ART']"

# --- 2. Widen column A ---------------------------------------------------
# NOTE: the COM width setter round-trips through a pixel-quantized font
# metric (Calibri 11) that adds ~5/6 of a character back on save, so we
# back that off here to land on the target stored width of exactly 25.
$ws.Columns.Item(1).ColumnWidth = 24.166666666666668

# --- 3. Append new rows 70-82 --------------------------------------------
$newRows = @(
    @("Media_repost_events",     "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Media CE: Reposting information about current events", "[0, 1]"),
    @("Media_sign_petition",     "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Media CE: Signing online petitions", "[0, 1]"),
    @("Media_follow_news",       "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Media CE: Following accounts/pages about current events", "[1, 0]"),
    @("Media_follow_officials",  "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Media CE: Following elected officials on social media", "[0, 1]"),
    @("Media_debate_opinions",   "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Media CE: Debating opinions with others who disagree", "[0, 1]"),
    @("Media_like_posts",        "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Media CE: Liking posts about current events", "[0, 1]"),
    @("Media_post_opinions",     "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Media CE: Posting own opinions about current events", "[0, 1]"),
    @("News_The_Elm",            "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Local news: The Elm", "[1, 0]"),
    @("News_Kent_County",        "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Local news: Kent County News", "[0, 1]"),
    @("News_Chestertown_Spy",    "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Local news: The Chestertown Spy", "[0, 1]"),
    @("News_Eastern_Shore_Post", "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Local news: Eastern Shore Post", "[0, 1]"),
    @("News_Shore_Daily",        "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Local news: Shore Daily News", "[0, 1]"),
    @("News_Delmarva_Now",       "—", "", "", "binary OHE", "1=selected, 0=not selected, -9=missing", -9, "Local news: Delmarva Now", "[0, 1]")
)

$startRow = 70
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
